$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the Heading1 title
#    paragraph. We clone the structure of the (soon to be removed) bold
#    "Play Dragon's Luck Power Reels for Free - Review" paragraph near the end
#    of the document via copy/paste so the leading empty run survives, then
#    retarget its text and append the rest of the description as a 3rd run.
# ---------------------------------------------------------------------------
$countBefore = $d.Paragraphs.Count
$pBoldSource = $d.Paragraphs($countBefore - 1)
$pBoldSource.Range.Copy()

$titlePar = $d.Paragraphs(1)
$titleEnd = $titlePar.Range
$titleEnd.Collapse(0)
$titleEnd.InsertParagraphAfter()

$metaPar = $d.Paragraphs(2)
$metaPar.Style = "Normal"
$metaRange = $metaPar.Range
$metaRange.Collapse(1)
$metaRange.Paste()

$metaPar = $d.Paragraphs(2)
[void]$metaPar.Range.Find.Execute("Play Dragon's Luck Power Reels for Free - Review", $true, $false, $false, $false, $false, $true, 1, $false, "Meta description", 2)

$metaPar = $d.Paragraphs(2)
$insertAt = $metaPar.Range.End - 1
$tail = $d.Range($insertAt, $insertAt)
[void]$tail.Font.Reset()
$tail.InsertAfter(": Read our review of Dragon's Luck Power Reels online slot game. Play for free and learn about gameplay, design, and payouts.")

# ---------------------------------------------------------------------------
# 2) Remove the duplicated bold "Play Dragon's Luck Power Reels for Free -
#    Review" paragraph that used to sit near the end of the document, right
#    before the italic meta-description paragraph.
# ---------------------------------------------------------------------------
$countNow = $d.Paragraphs.Count
$pBoldEnd = $d.Paragraphs($countNow - 1)
$pBoldEnd.Range.Delete()

# ---------------------------------------------------------------------------
# 3) Replace the final (italic) paragraph's text with the new image prompt.
#    Assign .Text directly (instead of Find/Replace) so Word's smart-quote
#    autocorrect doesn't mangle the apostrophes.
# ---------------------------------------------------------------------------
$countFinal = $d.Paragraphs.Count
$pPrompt = $d.Paragraphs($countFinal)
$promptRange = $pPrompt.Range
$promptTextRange = $d.Range($promptRange.Start, $promptRange.End - 1)
$promptTextRange.Text = "Prompt: Create a feature image for Dragon's Luck Power Reels in cartoon style featuring a happy Maya warrior with glasses. The image should incorporate the theme of Chinese tradition and the figure of the dragon in a visually appealing manner. It should also showcase the 10 reels and 30 paylines of the game, with the highest-valued symbol, the coin with the number 138, prominently displayed. The image should be optimized for use on both desktop computers and all iOS and Android mobile devices, and should capture the high volatility and distinctive gameplay features of the slot machine."
